# Arena Commits Table - add two new commit rows, extend the blank
# "buffer" rows below the table, and move the Total(h) row down to
# make room (commit: "creatures states update & Crafter NPC base code
# added").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Total" row (row 29) needs to end up at row 39, with nine
# freshly-formatted blank rows (29-37) inserted above it and row 38
# left completely empty. Inserting 10 rows at 29 (copying the
# formatting of row 28, the row above) achieves exactly that shift;
# we then wipe row 38 back to a truly blank (unformatted) row.
$ws.Rows("29:38").Insert()
$ws.Range("C38:G38").Clear()

# New commit log entries in the two rows that used to be blank
# placeholders (C26/C27 already carried the table's row styling).
$ws.Range("C26").Value = "game output system & levels update"
$ws.Range("G26").Value = 1.2

$ws.Range("C27").Value = "creatures states update & Crafter NPC base code"
$ws.Range("G27").Value = 2

# The Total(h) formula (now on row 39) needs to include the two new
# rows.
$ws.Range("G39").Formula = "=SUM(G4:G27)"

# Match the author's final on-screen selection/scroll position.
$ws.Range("H32").Select()
